$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quantity (F) and recompute value (G = rate * qty) for each changed item row ---
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = $ws.Range("D25").Value2 * $ws.Range("F25").Value2
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = $ws.Range("D28").Value2 * $ws.Range("F28").Value2
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $ws.Range("D30").Value2 * $ws.Range("F30").Value2
$ws.Range("F39").Value = 66
$ws.Range("G39").Value = $ws.Range("D39").Value2 * $ws.Range("F39").Value2
$ws.Range("F52").Value = 53
$ws.Range("G52").Value = $ws.Range("D52").Value2 * $ws.Range("F52").Value2
$ws.Range("F53").Value = 75
$ws.Range("G53").Value = $ws.Range("D53").Value2 * $ws.Range("F53").Value2
$ws.Range("F86").Value = 65
$ws.Range("G86").Value = $ws.Range("D86").Value2 * $ws.Range("F86").Value2
$ws.Range("F94").Value = 126
$ws.Range("G94").Value = $ws.Range("D94").Value2 * $ws.Range("F94").Value2
$ws.Range("F105").Value = 187
$ws.Range("G105").Value = $ws.Range("D105").Value2 * $ws.Range("F105").Value2
$ws.Range("F110").Value = 98
$ws.Range("G110").Value = $ws.Range("D110").Value2 * $ws.Range("F110").Value2
$ws.Range("F133").Value = 38
$ws.Range("G133").Value = $ws.Range("D133").Value2 * $ws.Range("F133").Value2
$ws.Range("F148").Value = 32
$ws.Range("G148").Value = $ws.Range("D148").Value2 * $ws.Range("F148").Value2
$ws.Range("F154").Value = 9
$ws.Range("G154").Value = $ws.Range("D154").Value2 * $ws.Range("F154").Value2
$ws.Range("F157").Value = 113
$ws.Range("G157").Value = $ws.Range("D157").Value2 * $ws.Range("F157").Value2
$ws.Range("F173").Value = 6
$ws.Range("G173").Value = $ws.Range("D173").Value2 * $ws.Range("F173").Value2
$ws.Range("F186").Value = 10
$ws.Range("G186").Value = $ws.Range("D186").Value2 * $ws.Range("F186").Value2
$ws.Range("F191").Value = 5
$ws.Range("G191").Value = $ws.Range("D191").Value2 * $ws.Range("F191").Value2
$ws.Range("F200").Value = 8
$ws.Range("G200").Value = $ws.Range("D200").Value2 * $ws.Range("F200").Value2
$ws.Range("F201").Value = 66
$ws.Range("G201").Value = $ws.Range("D201").Value2 * $ws.Range("F201").Value2
$ws.Range("F211").Value = 23
$ws.Range("G211").Value = $ws.Range("D211").Value2 * $ws.Range("F211").Value2
$ws.Range("F219").Value = 64
$ws.Range("G219").Value = $ws.Range("D219").Value2 * $ws.Range("F219").Value2
$ws.Range("F236").Value = 0
$ws.Range("G236").Value = $ws.Range("D236").Value2 * $ws.Range("F236").Value2
$ws.Range("F239").Value = 35
$ws.Range("G239").Value = $ws.Range("D239").Value2 * $ws.Range("F239").Value2
$ws.Range("F243").Value = 3
$ws.Range("G243").Value = $ws.Range("D243").Value2 * $ws.Range("F243").Value2
$ws.Range("F244").Value = 21
$ws.Range("G244").Value = $ws.Range("D244").Value2 * $ws.Range("F244").Value2
$ws.Range("F264").Value = 173
$ws.Range("G264").Value = $ws.Range("D264").Value2 * $ws.Range("F264").Value2
$ws.Range("F303").Value = 67
$ws.Range("G303").Value = $ws.Range("D303").Value2 * $ws.Range("F303").Value2
$ws.Range("F304").Value = 29
$ws.Range("G304").Value = $ws.Range("D304").Value2 * $ws.Range("F304").Value2
$ws.Range("F307").Value = 21
$ws.Range("G307").Value = $ws.Range("D307").Value2 * $ws.Range("F307").Value2
$ws.Range("F310").Value = 177
$ws.Range("G310").Value = $ws.Range("D310").Value2 * $ws.Range("F310").Value2
$ws.Range("F319").Value = 133
$ws.Range("G319").Value = $ws.Range("D319").Value2 * $ws.Range("F319").Value2
$ws.Range("F328").Value = 167
$ws.Range("G328").Value = $ws.Range("D328").Value2 * $ws.Range("F328").Value2
$ws.Range("F340").Value = 65
$ws.Range("G340").Value = $ws.Range("D340").Value2 * $ws.Range("F340").Value2
$ws.Range("F342").Value = 30
$ws.Range("G342").Value = $ws.Range("D342").Value2 * $ws.Range("F342").Value2
$ws.Range("F347").Value = 18
$ws.Range("G347").Value = $ws.Range("D347").Value2 * $ws.Range("F347").Value2
$ws.Range("F356").Value = 28
$ws.Range("G356").Value = $ws.Range("D356").Value2 * $ws.Range("F356").Value2
$ws.Range("F391").Value = 46
$ws.Range("G391").Value = $ws.Range("D391").Value2 * $ws.Range("F391").Value2
$ws.Range("F412").Value = 77
$ws.Range("G412").Value = $ws.Range("D412").Value2 * $ws.Range("F412").Value2
$ws.Range("F420").Value = 43
$ws.Range("G420").Value = $ws.Range("D420").Value2 * $ws.Range("F420").Value2
$ws.Range("F421").Value = 43
$ws.Range("G421").Value = $ws.Range("D421").Value2 * $ws.Range("F421").Value2
$ws.Range("F430").Value = 18
$ws.Range("G430").Value = $ws.Range("D430").Value2 * $ws.Range("F430").Value2
$ws.Range("F435").Value = 79
$ws.Range("G435").Value = $ws.Range("D435").Value2 * $ws.Range("F435").Value2
$ws.Range("F458").Value = 512
$ws.Range("G458").Value = $ws.Range("D458").Value2 * $ws.Range("F458").Value2
$ws.Range("F459").Value = 488
$ws.Range("G459").Value = $ws.Range("D459").Value2 * $ws.Range("F459").Value2
$ws.Range("F460").Value = 569
$ws.Range("G460").Value = $ws.Range("D460").Value2 * $ws.Range("F460").Value2
$ws.Range("F461").Value = 272
$ws.Range("G461").Value = $ws.Range("D461").Value2 * $ws.Range("F461").Value2
$ws.Range("F463").Value = 323
$ws.Range("G463").Value = $ws.Range("D463").Value2 * $ws.Range("F463").Value2
$ws.Range("F465").Value = 353
$ws.Range("G465").Value = $ws.Range("D465").Value2 * $ws.Range("F465").Value2
$ws.Range("F466").Value = 376
$ws.Range("G466").Value = $ws.Range("D466").Value2 * $ws.Range("F466").Value2
$ws.Range("F467").Value = 443
$ws.Range("G467").Value = $ws.Range("D467").Value2 * $ws.Range("F467").Value2
$ws.Range("F468").Value = 101
$ws.Range("G468").Value = $ws.Range("D468").Value2 * $ws.Range("F468").Value2
$ws.Range("F469").Value = 970
$ws.Range("G469").Value = $ws.Range("D469").Value2 * $ws.Range("F469").Value2
$ws.Range("F474").Value = 274
$ws.Range("G474").Value = $ws.Range("D474").Value2 * $ws.Range("F474").Value2
$ws.Range("F475").Value = 578
$ws.Range("G475").Value = $ws.Range("D475").Value2 * $ws.Range("F475").Value2
$ws.Range("F519").Value = 577
$ws.Range("G519").Value = $ws.Range("D519").Value2 * $ws.Range("F519").Value2
$ws.Range("F521").Value = 351
$ws.Range("G521").Value = $ws.Range("D521").Value2 * $ws.Range("F521").Value2
$ws.Range("F526").Value = 140
$ws.Range("G526").Value = $ws.Range("D526").Value2 * $ws.Range("F526").Value2
$ws.Range("F538").Value = 16
$ws.Range("G538").Value = $ws.Range("D538").Value2 * $ws.Range("F538").Value2
$ws.Range("F539").Value = 23
$ws.Range("G539").Value = $ws.Range("D539").Value2 * $ws.Range("F539").Value2
$ws.Range("F556").Value = 55
$ws.Range("G556").Value = $ws.Range("D556").Value2 * $ws.Range("F556").Value2
$ws.Range("F557").Value = 103
$ws.Range("G557").Value = $ws.Range("D557").Value2 * $ws.Range("F557").Value2
$ws.Range("F559").Value = 33
$ws.Range("G559").Value = $ws.Range("D559").Value2 * $ws.Range("F559").Value2
$ws.Range("F560").Value = 18
$ws.Range("G560").Value = $ws.Range("D560").Value2 * $ws.Range("F560").Value2
$ws.Range("F572").Value = 14
$ws.Range("G572").Value = $ws.Range("D572").Value2 * $ws.Range("F572").Value2
$ws.Range("F580").Value = 24
$ws.Range("G580").Value = $ws.Range("D580").Value2 * $ws.Range("F580").Value2
$ws.Range("F597").Value = 75
$ws.Range("G597").Value = $ws.Range("D597").Value2 * $ws.Range("F597").Value2
$ws.Range("F599").Value = 91
$ws.Range("G599").Value = $ws.Range("D599").Value2 * $ws.Range("F599").Value2
$ws.Range("F671").Value = 4
$ws.Range("G671").Value = $ws.Range("D671").Value2 * $ws.Range("F671").Value2
$ws.Range("F689").Value = 353
$ws.Range("G689").Value = $ws.Range("D689").Value2 * $ws.Range("F689").Value2
$ws.Range("F690").Value = 137
$ws.Range("G690").Value = $ws.Range("D690").Value2 * $ws.Range("F690").Value2
$ws.Range("F692").Value = 17
$ws.Range("G692").Value = $ws.Range("D692").Value2 * $ws.Range("F692").Value2
$ws.Range("F693").Value = 364
$ws.Range("G693").Value = $ws.Range("D693").Value2 * $ws.Range("F693").Value2
$ws.Range("F694").Value = 20
$ws.Range("G694").Value = $ws.Range("D694").Value2 * $ws.Range("F694").Value2
$ws.Range("F698").Value = 176
$ws.Range("G698").Value = $ws.Range("D698").Value2 * $ws.Range("F698").Value2
$ws.Range("F699").Value = 383
$ws.Range("G699").Value = $ws.Range("D699").Value2 * $ws.Range("F699").Value2
$ws.Range("F701").Value = 197
$ws.Range("G701").Value = $ws.Range("D701").Value2 * $ws.Range("F701").Value2
$ws.Range("F702").Value = 13
$ws.Range("G702").Value = $ws.Range("D702").Value2 * $ws.Range("F702").Value2
$ws.Range("F706").Value = 707
$ws.Range("G706").Value = $ws.Range("D706").Value2 * $ws.Range("F706").Value2
$ws.Range("F707").Value = 65
$ws.Range("G707").Value = $ws.Range("D707").Value2 * $ws.Range("F707").Value2
$ws.Range("F708").Value = 785
$ws.Range("G708").Value = $ws.Range("D708").Value2 * $ws.Range("F708").Value2
$ws.Range("F719").Value = 204
$ws.Range("G719").Value = $ws.Range("D719").Value2 * $ws.Range("F719").Value2
$ws.Range("F726").Value = 72
$ws.Range("G726").Value = $ws.Range("D726").Value2 * $ws.Range("F726").Value2
$ws.Range("F731").Value = 9
$ws.Range("G731").Value = $ws.Range("D731").Value2 * $ws.Range("F731").Value2
$ws.Range("F733").Value = 28
$ws.Range("G733").Value = $ws.Range("D733").Value2 * $ws.Range("F733").Value2
$ws.Range("F735").Value = 54
$ws.Range("G735").Value = $ws.Range("D735").Value2 * $ws.Range("F735").Value2
$ws.Range("F771").Value = 1
$ws.Range("G771").Value = $ws.Range("D771").Value2 * $ws.Range("F771").Value2
$ws.Range("F775").Value = 2
$ws.Range("G775").Value = $ws.Range("D775").Value2 * $ws.Range("F775").Value2
$ws.Range("F781").Value = 503
$ws.Range("G781").Value = $ws.Range("D781").Value2 * $ws.Range("F781").Value2
$ws.Range("F785").Value = 11
$ws.Range("G785").Value = $ws.Range("D785").Value2 * $ws.Range("F785").Value2

# --- Rows 300 and 301 had their data (code/rate/MRP/qty/value) swapped between each other ---
$cols = @("B","D","E","F","G")
foreach ($col in $cols) {
    $ref300 = "${col}300"
    $ref301 = "${col}301"
    $tmp = $ws.Range($ref300).Value2
    $ws.Range($ref300).Value = $ws.Range($ref301).Value2
    $ws.Range($ref301).Value = $tmp
}

$sumFn = $excel.WorksheetFunction

# --- Recompute company "Sub Total:" rows affected by the quantity changes above ---
$ws.Range("B34").Value = $sumFn.Sum($ws.Range("G21:G33"))
$ws.Range("B63").Value = $sumFn.Sum($ws.Range("G36:G62"))
$ws.Range("B122").Value = $sumFn.Sum($ws.Range("G84:G121"))
$ws.Range("B140").Value = $sumFn.Sum($ws.Range("G131:G139"))
$ws.Range("B160").Value = $sumFn.Sum($ws.Range("G148:G159"))
$ws.Range("B181").Value = $sumFn.Sum($ws.Range("G170:G180"))
$ws.Range("B202").Value = $sumFn.Sum($ws.Range("G183:G201"))
$ws.Range("B214").Value = $sumFn.Sum($ws.Range("G204:G213"))
$ws.Range("B221").Value = $sumFn.Sum($ws.Range("G216:G220"))
$ws.Range("B245").Value = $sumFn.Sum($ws.Range("G236:G244"))
$ws.Range("B296").Value = $sumFn.Sum($ws.Range("G257:G295"))
$ws.Range("B366").Value = $sumFn.Sum($ws.Range("G298:G365"))
$ws.Range("B397").Value = $sumFn.Sum($ws.Range("G391:G396"))
$ws.Range("B424").Value = $sumFn.Sum($ws.Range("G411:G423"))
$ws.Range("B441").Value = $sumFn.Sum($ws.Range("G426:G440"))
$ws.Range("B476").Value = $sumFn.Sum($ws.Range("G458:G475"))
$ws.Range("B527").Value = $sumFn.Sum($ws.Range("G519:G526"))
$ws.Range("B548").Value = $sumFn.Sum($ws.Range("G534:G547"))
$ws.Range("B565").Value = $sumFn.Sum($ws.Range("G553:G564"))
$ws.Range("B577").Value = $sumFn.Sum($ws.Range("G567:G576"))
$ws.Range("B594").Value = $sumFn.Sum($ws.Range("G579:G593"))
$ws.Range("B603").Value = $sumFn.Sum($ws.Range("G596:G602"))
$ws.Range("B674").Value = $sumFn.Sum($ws.Range("G667:G673"))
$ws.Range("B710").Value = $sumFn.Sum($ws.Range("G689:G709"))
$ws.Range("B728").Value = $sumFn.Sum($ws.Range("G712:G727"))
$ws.Range("B736").Value = $sumFn.Sum($ws.Range("G730:G735"))
$ws.Range("B772").Value = $sumFn.Sum($ws.Range("G750:G771"))
$ws.Range("B779").Value = $sumFn.Sum($ws.Range("G774:G778"))
$ws.Range("B786").Value = $sumFn.Sum($ws.Range("G781:G785"))

# --- Recompute the overall "Sub Total:" (row 792) and "Grand Total:" (row 793) ---
$subtotalCells = @("B7","B12","B19","B34","B63","B66","B70","B73","B82","B122","B126","B129","B140","B146","B160","B168","B181","B202","B214","B221","B224","B229","B234","B245","B255","B296","B366","B375","B378","B389","B397","B406","B409","B424","B441","B450","B456","B476","B483","B501","B505","B514","B517","B527","B532","B548","B551","B565","B577","B594","B603","B620","B631","B637","B650","B653","B658","B665","B674","B678","B687","B710","B728","B736","B741","B748","B772","B779","B786","B791")
$grand = 0
foreach ($cellRef in $subtotalCells) {
    $grand = $grand + $ws.Range($cellRef).Value2
}
$ws.Range("B792").Value = $grand
$ws.Range("B793").Value = $grand